$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.547.62'
$ws.Range("E2").Value = '  -6.23%  '
$ws.Range("D3").Value = '1.809.00'
$ws.Range("E3").Value = '  -5.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '276.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -9.89%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -6.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3494'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -8.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.97'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06634'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -9.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.97'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -10.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8346'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -7.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07816'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.72%  '
$ws.Range("D14").Value = '1.805.53'
$ws.Range("E14").Value = '  +70.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.030'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.01'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -9.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.85'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.01%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007848'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -9.32%  '
$ws.Range("D21").Value = '25.612.55'
$ws.Range("E21").Value = '  -6.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.714'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.65%  '
$ws.Range("D23").Value = '2.038.63'
$ws.Range("E23").Value = '  +67.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.933'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -7.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.055'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.95'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.667'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.097'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -8.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.88'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '108.71'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -6.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.283'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -11.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.197'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -11.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08806'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04789'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7313'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -11.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.868'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.117'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -8.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.001'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.034'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -8.92%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01856'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.33%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5193'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -11.91%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.309'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -13.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9562'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -11.37%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '110.98'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.04%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.158'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.038'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -14.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4587'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1377'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -9.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.237'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.52'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.49%  '
